$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) for two events, on both the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet,
# mirroring the same rows: F2 535 -> 536, F7 751 -> 752.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 536
    $ws.Range("F7").Value = 752
}
